$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.472.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.62%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.945.88'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '508.15'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.76'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.735'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.31%  '

$ws.Range("E10").Value = '  +4.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000350'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.66%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.69'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.77%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.572.53'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.33%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.957.52'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.04%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.24%  '

$ws.Range("E17").Value = '  -0.26%  '

$ws.Range("E18").Value = '  +7.26%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.38%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.503.10'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '436.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.72'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.14'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.07%  '

$ws.Range("E26").Value = '  +6.94%  '

$ws.Range("E27").Value = '  -2.42%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.61%  '

$ws.Range("E29").Value = '  -2.80%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '707.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.94%  '

$ws.Range("E32").Value = '  -0.85%  '

$ws.Range("E33").Value = '  -0.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '65.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.450'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +13.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0881'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.09'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '40.93'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.78%  '

$ws.Range("E39").Value = '  +0.80%  '

$ws.Range("E40").Value = '  -0.12%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("E42").Value = '  +2.00%  '

$ws.Range("E43").Value = '  -1.85%  '

$ws.Range("E44").Value = '  +5.33%  '

$ws.Range("E45").Value = '  -4.53%  '

$ws.Range("E46").Value = '  +1.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.38'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.64%  '

$ws.Range("E48").Value = '  +6.13%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0353'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.75%  '

$ws.Range("E51").Value = '  -1.23%  '
